$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 ("Marking"): update Right and Wrong marking weights
$ws.Range("B11").Value = 9
$ws.Range("C11").Value = 2

# Row 12 ("Total"): recalculated totals based on the new marking weights
$ws.Range("B12").Value = 216
$ws.Range("C12").Value = -8
$ws.Range("E12").Value = "208/252"
